# Add two new columns, I ("I0") and J ("IF"), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - style matches the other header cells (B1:H1), i.e. the
# "Header" cell style: bold font, thin border, centered/top aligned.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Data values for I2:I67 and J2:J67 (rows 2-67).
$I = @(7,9,8,6,8,8,9,8,5,12,2,8,7,5,9,7,6,9,9,8,6,7,7,8,8,5,8,6,5,10,7,5,7,7,9,9,7,7,5,9,9,7,8,11,6,9,8,9,7,7,7,6,9,5,6,7,5,7,6,9,6,8,8,8,7,7)
$J = @(8,9,8,6,8,8,9,8,5,12,2,8,7,5,9,7,6,9,9,8,6,7,7,8,8,5,8,6,5,10,7,5,7,7,9,9,7,7,5,9,9,7,8,11,6,9,8,9,7,7,7,6,9,5,6,7,5,7,6,9,6,8,8,8,7,7)

for ($r = 2; $r -le 67; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $I[$idx]
    $ws.Cells.Item($r, 10).Value = $J[$idx]
}
